$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three trailing rows whose sending cluster was "ECs" (rows 8-10)
$ws.Rows("8:10").Delete()

# Update rows 2-7 with the refreshed TPM-derived values
# Row 2
$ws.Cells.Item(2,1).Value2 = "FAPs"
$ws.Cells.Item(2,2).Value2 = "Slit1"
$ws.Cells.Item(2,3).Value2 = "Robo1"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 2
$ws.Cells.Item(2,6).Value2 = 0.6666666666666666
$ws.Cells.Item(2,7).Value2 = 0.2055996666666667
$ws.Cells.Item(2,8).Value2 = 0.616799
$ws.Cells.Item(2,9).Value2 = 0.9059768423248155
$ws.Cells.Item(2,10).Value2 = 0.9059768423248156
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 0.040495
$ws.Cells.Item(2,14).Value2 = 0.121485
$ws.Cells.Item(2,15).Value2 = 0.002191743187342868
$ws.Cells.Item(2,16).Value2 = 0.002191743187342869
$ws.Cells.Item(2,17).Value2 = 0.008325758501666666
$ws.Cells.Item(2,18).Value2 = 0.074931826515
$ws.Cells.Item(2,19).Value2 = 0.001985668572055818
$ws.Cells.Item(2,20).Value2 = 0.001985668572055819

# Row 3
$ws.Cells.Item(3,1).Value2 = "FAPs"
$ws.Cells.Item(3,2).Value2 = "Slit1"
$ws.Cells.Item(3,3).Value2 = "Robo1"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 2
$ws.Cells.Item(3,6).Value2 = 0.6666666666666666
$ws.Cells.Item(3,7).Value2 = 0.2055996666666667
$ws.Cells.Item(3,8).Value2 = 0.616799
$ws.Cells.Item(3,9).Value2 = 0.9059768423248155
$ws.Cells.Item(3,10).Value2 = 0.9059768423248156
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 16.98312366666667
$ws.Cells.Item(3,14).Value2 = 50.949371
$ws.Cells.Item(3,15).Value2 = 0.9191911494312409
$ws.Cells.Item(3,16).Value2 = 0.9191911494312409
$ws.Cells.Item(3,17).Value2 = 3.491724564825444
$ws.Cells.Item(3,18).Value2 = 31.425521083429
$ws.Cells.Item(3,19).Value2 = 0.8327658950546333
$ws.Cells.Item(3,20).Value2 = 0.8327658950546334

# Row 4
$ws.Cells.Item(4,1).Value2 = "FAPs"
$ws.Cells.Item(4,2).Value2 = "Slit1"
$ws.Cells.Item(4,3).Value2 = "Robo1"
$ws.Cells.Item(4,4).Value2 = "MuSCs"
$ws.Cells.Item(4,5).Value2 = 2
$ws.Cells.Item(4,6).Value2 = 0.6666666666666666
$ws.Cells.Item(4,7).Value2 = 0.2055996666666667
$ws.Cells.Item(4,8).Value2 = 0.616799
$ws.Cells.Item(4,9).Value2 = 0.9059768423248155
$ws.Cells.Item(4,10).Value2 = 0.9059768423248156
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 1.452542333333333
$ws.Cells.Item(4,14).Value2 = 4.357627
$ws.Cells.Item(4,15).Value2 = 0.07861710738141615
$ws.Cells.Item(4,16).Value2 = 0.07861710738141615
$ws.Cells.Item(4,17).Value2 = 0.2986422195525555
$ws.Cells.Item(4,18).Value2 = 2.687779975973
$ws.Cells.Item(4,19).Value2 = 0.07122527869812635
$ws.Cells.Item(4,20).Value2 = 0.07122527869812635

# Row 5
$ws.Cells.Item(5,1).Value2 = "MuSCs"
$ws.Cells.Item(5,2).Value2 = "Slit1"
$ws.Cells.Item(5,3).Value2 = "Robo1"
$ws.Cells.Item(5,4).Value2 = "ECs"
$ws.Cells.Item(5,5).Value2 = 2
$ws.Cells.Item(5,6).Value2 = 0.6666666666666666
$ws.Cells.Item(5,7).Value2 = 0.02133733333333333
$ws.Cells.Item(5,8).Value2 = 0.064012
$ws.Cells.Item(5,9).Value2 = 0.09402315767518445
$ws.Cells.Item(5,10).Value2 = 0.09402315767518446
$ws.Cells.Item(5,11).Value2 = 3
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 0.040495
$ws.Cells.Item(5,14).Value2 = 0.121485
$ws.Cells.Item(5,15).Value2 = 0.002191743187342868
$ws.Cells.Item(5,16).Value2 = 0.002191743187342869
$ws.Cells.Item(5,17).Value2 = 0.0008640553133333332
$ws.Cells.Item(5,18).Value2 = 0.00777649782
$ws.Cells.Item(5,19).Value2 = 0.0002060746152870498
$ws.Cells.Item(5,20).Value2 = 0.0002060746152870499

# Row 6
$ws.Cells.Item(6,1).Value2 = "MuSCs"
$ws.Cells.Item(6,2).Value2 = "Slit1"
$ws.Cells.Item(6,3).Value2 = "Robo1"
$ws.Cells.Item(6,4).Value2 = "FAPs"
$ws.Cells.Item(6,5).Value2 = 2
$ws.Cells.Item(6,6).Value2 = 0.6666666666666666
$ws.Cells.Item(6,7).Value2 = 0.02133733333333333
$ws.Cells.Item(6,8).Value2 = 0.064012
$ws.Cells.Item(6,9).Value2 = 0.09402315767518445
$ws.Cells.Item(6,10).Value2 = 0.09402315767518446
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 16.98312366666667
$ws.Cells.Item(6,14).Value2 = 50.949371
$ws.Cells.Item(6,15).Value2 = 0.9191911494312409
$ws.Cells.Item(6,16).Value2 = 0.9191911494312409
$ws.Cells.Item(6,17).Value2 = 0.3623745707168889
$ws.Cells.Item(6,18).Value2 = 3.261371136452
$ws.Cells.Item(6,19).Value2 = 0.0864252543766076
$ws.Cells.Item(6,20).Value2 = 0.0864252543766076

# Row 7
$ws.Cells.Item(7,1).Value2 = "MuSCs"
$ws.Cells.Item(7,2).Value2 = "Slit1"
$ws.Cells.Item(7,3).Value2 = "Robo1"
$ws.Cells.Item(7,4).Value2 = "MuSCs"
$ws.Cells.Item(7,5).Value2 = 2
$ws.Cells.Item(7,6).Value2 = 0.6666666666666666
$ws.Cells.Item(7,7).Value2 = 0.02133733333333333
$ws.Cells.Item(7,8).Value2 = 0.064012
$ws.Cells.Item(7,9).Value2 = 0.09402315767518445
$ws.Cells.Item(7,10).Value2 = 0.09402315767518446
$ws.Cells.Item(7,11).Value2 = 3
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 1.452542333333333
$ws.Cells.Item(7,14).Value2 = 4.357627
$ws.Cells.Item(7,15).Value2 = 0.07861710738141615
$ws.Cells.Item(7,16).Value2 = 0.07861710738141615
$ws.Cells.Item(7,17).Value2 = 0.03099337994711111
$ws.Cells.Item(7,18).Value2 = 0.278940419524
$ws.Cells.Item(7,19).Value2 = 0.007391828683289798
$ws.Cells.Item(7,20).Value2 = 0.007391828683289799

